$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

# Sheet2: rename the group (质控组 -> 北京组) and update the K2 value
$ws2.Range("A2").Value = "北京组"
$ws2.Range("K2").Value = 5.87

# Sheet1: remove 张悦's row entirely (row 2); remaining rows shift up
$ws1.Rows(2).Delete()

# Sheet1: rename the group (质控组 -> 北京组) for the remaining rows
$ws1.Range("A2:A3").Value = "北京组"

# Restore the view/selection state: Sheet2's selection first (it is the
# default-active sheet), then Sheet1's - selecting on Sheet1 makes it the
# active tab, matching the target workbook state.
$ws2.Range("B25").Select() | Out-Null
$ws1.Range("B4").Select() | Out-Null
